$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44285
$ws.Range("L2").Value = 'Segunda'
$ws.Range("M2").Value = 160
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15500
$ws.Range("S2").Value = 1033

$ws.Range("D3").Value = 44348
$ws.Range("L3").Value = 'Primera'
$ws.Range("M3").Value = 100
$ws.Range("N3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 1000

$ws.Range("L4").Value = 'Segunda'
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 13000
$ws.Range("O4").Value = 13500
$ws.Range("P4").Value = 13250
$ws.Range("S4").Value = 883

$ws.Range("D5").Value = 44299
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 60
$ws.Range("N5").Value = 17500
$ws.Range("O5").Value = 17500
$ws.Range("P5").Value = 17500
$ws.Range("S5").Value = 1167

$ws.Range("D6").Value = 44299
$ws.Range("M6").Value = 120
$ws.Range("N6").Value = 14000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 14500
$ws.Range("S6").Value = 967

$ws.Range("D7").Value = 44344
$ws.Range("L7").Value = 'Primera'
$ws.Range("N7").Value = 16000
$ws.Range("O7").Value = 16000
$ws.Range("P7").Value = 16000
$ws.Range("Q7").Value = '$/caja 15 kilos empedrada'
$ws.Range("S7").Value = 1067
$ws.Range("T7").Value = 15

$ws.Range("D8").Value = 44344
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 120
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13500
$ws.Range("P8").Value = 13250
$ws.Range("Q8").Value = '$/caja 15 kilos empedrada'
$ws.Range("S8").Value = 883
$ws.Range("T8").Value = 15

$ws.Range("L9").Value = 'Especial'
$ws.Range("M9").Value = 100
$ws.Range("N9").Value = 17500
$ws.Range("O9").Value = 17500
$ws.Range("P9").Value = 17500
$ws.Range("S9").Value = 1250

$ws.Range("D10").Value = 44313
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 16000
$ws.Range("P10").Value = 16000
$ws.Range("Q10").Value = '$/caja 14 kilos empedrada'
$ws.Range("S10").Value = 1143
$ws.Range("T10").Value = 14

$ws.Range("D11").Value = 44313
$ws.Range("M11").Value = 80
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("Q11").Value = '$/caja 14 kilos empedrada'
$ws.Range("S11").Value = 1000
$ws.Range("T11").Value = 14

$ws.Range("D12").Value = 44327
$ws.Range("N12").Value = 17000
$ws.Range("O12").Value = 17000
$ws.Range("P12").Value = 17000
$ws.Range("S12").Value = 1133

$ws.Range("D13").Value = 44327

$ws.Range("D14").Value = 44336
$ws.Range("M14").Value = 60

$ws.Range("D15").Value = 44336
$ws.Range("M15").Value = 120
$ws.Range("O15").Value = 14500
$ws.Range("P15").Value = 14250
$ws.Range("S15").Value = 950

$ws.Range("D16").Value = 44309
$ws.Range("M16").Value = 100

$ws.Range("D17").Value = 44309
$ws.Range("M17").Value = 200
$ws.Range("O17").Value = 14500
$ws.Range("P17").Value = 14250
$ws.Range("S17").Value = 950

$ws.Range("D18").Value = 44298
$ws.Range("L18").Value = 'Segunda'
$ws.Range("M18").Value = 80

$ws.Range("D19").Value = 44305
$ws.Range("M19").Value = 60
$ws.Range("N19").Value = 17500
$ws.Range("O19").Value = 17500
$ws.Range("P19").Value = 17500
$ws.Range("S19").Value = 1167

$ws.Range("D20").Value = 44305
$ws.Range("M20").Value = 120
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("S20").Value = 967

$ws.Range("D21").Value = 44334
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 14000
$ws.Range("P21").Value = 15500
$ws.Range("S21").Value = 1033

$ws.Range("D22").Value = 44334
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 14500
$ws.Range("P22").Value = 14500
$ws.Range("S22").Value = 967

$ws.Range("D23").Value = 44316
$ws.Range("N23").Value = 17500
$ws.Range("O23").Value = 17500
$ws.Range("P23").Value = 17500
$ws.Range("S23").Value = 1167

$ws.Range("D24").Value = 44316

$ws.Range("D25").Value = 44351
$ws.Range("M25").Value = 100
$ws.Range("N25").Value = 15000
$ws.Range("P25").Value = 15000
$ws.Range("S25").Value = 1000

$ws.Range("D26").Value = 44351
$ws.Range("L26").Value = 'Segunda'
$ws.Range("M26").Value = 200
$ws.Range("N26").Value = 13000
$ws.Range("O26").Value = 13500
$ws.Range("P26").Value = 13250
$ws.Range("S26").Value = 883

$ws.Range("D27").Value = 44293
$ws.Range("L27").Value = 'Primera'
$ws.Range("M27").Value = 60
$ws.Range("O27").Value = 15000
$ws.Range("P27").Value = 14500
$ws.Range("S27").Value = 967

$ws.Range("D29").Value = 44323
$ws.Range("N29").Value = 17000
$ws.Range("O29").Value = 17000
$ws.Range("P29").Value = 17000
$ws.Range("S29").Value = 1133

$ws.Range("D30").Value = 44323
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 14000
$ws.Range("O30").Value = 14000
$ws.Range("P30").Value = 14000
$ws.Range("S30").Value = 933

$ws.Range("D31").Value = 44306
$ws.Range("N31").Value = 17500
$ws.Range("O31").Value = 17500
$ws.Range("P31").Value = 17500
$ws.Range("S31").Value = 1167

$ws.Range("D32").Value = 44306
$ws.Range("M32").Value = 200
$ws.Range("N32").Value = 14000
$ws.Range("O32").Value = 14500
$ws.Range("P32").Value = 14250
$ws.Range("S32").Value = 950

$ws.Range("D33").Value = 44301
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 17500
$ws.Range("O33").Value = 17500
$ws.Range("P33").Value = 17500
$ws.Range("S33").Value = 1167

$ws.Range("D34").Value = 44301
$ws.Range("M34").Value = 80
$ws.Range("N34").Value = 14000
$ws.Range("O34").Value = 15000

$ws.Range("D35").Value = 44330
$ws.Range("N35").Value = 17000
$ws.Range("O35").Value = 17000
$ws.Range("P35").Value = 17000
$ws.Range("S35").Value = 1133

$ws.Range("D36").Value = 44330

$ws.Range("D37").Value = 44295
$ws.Range("M37").Value = 160
$ws.Range("N37").Value = 14000
$ws.Range("O37").Value = 15000
$ws.Range("P37").Value = 14500
$ws.Range("S37").Value = 967

$ws.Range("D38").Value = 44302
$ws.Range("L38").Value = 'Primera'
$ws.Range("M38").Value = 100
$ws.Range("N38").Value = 17500
$ws.Range("O38").Value = 17500
$ws.Range("P38").Value = 17500
$ws.Range("S38").Value = 1167

$ws.Range("D39").Value = 44302
$ws.Range("M39").Value = 200
